$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "2016 Grant Awards By Agency Table"

$ws.Range("F4").Select() | Out-Null
